# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held the literal sheet/file label "6-18-2007-08"
# for every data row; correct it to the actual ISO game date 2008-06-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF
    if ($cell.Value() -eq "6-18-2007-08") {
        # Force text storage so Excel doesn't reinterpret the
        # "yyyy-mm-dd"-looking string as a date serial number, then
        # restore the default "Normal" style so no stray cell style
        # (numFmt) is left behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = "2008-06-18"
        $cell.Style = "Normal"
    }
}
